# Denmark Division 1 - league base update (17-02-2024 22:47)
#
# The underlying change re-sorts a handful of same-kickoff-time match
# records: for several rounds, two (or three) fixtures played at the same
# Date/Div were re-ordered, so row N must now hold the data that used to
# live in row M (id/date/div columns A-E stay put; B and F..AC - the
# match id, teams, score, result and all odds columns - move together).
#
# We implement this by swapping (or cyclically rotating) the B,F:AC
# ranges of each affected row group directly on the live sheet, so the
# values always come from the workbook itself rather than being
# hard-coded twice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowRecord($row) {
    # Columns B (2) and F..AC (6..29) hold the match id / teams / score /
    # odds payload that travels together when rows are re-ordered.
    $rec = @()
    $rec += $ws.Cells.Item($row, 2).Value2
    for ($c = 6; $c -le 29; $c++) {
        $rec += $ws.Cells.Item($row, $c).Value2
    }
    return $rec
}

function Set-RowRecord($row, $rec) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $idx = 1
    for ($c = 6; $c -le 29; $c++) {
        $ws.Cells.Item($row, $c).Value = $rec[$idx]
        $idx++
    }
}

function Swap-Rows($r1, $r2) {
    $rec1 = Get-RowRecord $r1
    $rec2 = Get-RowRecord $r2
    Set-RowRecord $r1 $rec2
    Set-RowRecord $r2 $rec1
}

function Rotate-Rows($rows) {
    # rows[0] <- rows[1] <- ... <- rows[n-1] <- rows[0]
    $n = $rows.Length
    $recs = @()
    for ($i = 0; $i -lt $n; $i++) {
        $recs += ,(Get-RowRecord $rows[$i])
    }
    for ($i = 0; $i -lt $n; $i++) {
        $srcIdx = ($i + 1) % $n
        Set-RowRecord $rows[$i] $recs[$srcIdx]
    }
}

# Simple pairwise swaps
Swap-Rows 9 10
Swap-Rows 31 32
Swap-Rows 70 71
Swap-Rows 84 85
Swap-Rows 87 88
Swap-Rows 122 123
Swap-Rows 125 126
Swap-Rows 142 143
Swap-Rows 152 153
Swap-Rows 159 160
Swap-Rows 166 167
Swap-Rows 198 199

# Three-way cyclic rotations
Rotate-Rows @(89, 90, 91)
Rotate-Rows @(146, 147, 148)
